$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 10-14 (the scheduling re-run covers fewer, more consolidated rows)
$ws.Rows("10:14").Delete()

# Update rows 4-9 with the recomputed scheduling data (more zones covered per the commit message)
$ws.Range("A4").Value = 251346
$ws.Range("B4").Value = "BIMEC 2"
$ws.Range("C4").Value = 36.5
$ws.Range("D4").Value = 70.16363636363636
$ws.Range("E4").Value = "2025-04-10 10:41:00"
$ws.Range("F4").Value = "2025-04-10 11:17:30"
$ws.Range("G4").Value = "2025-04-10 11:17:30"
$ws.Range("H4").Value = "2025-04-10 12:27:39"
$ws.Range("I4").Value = 3859
$ws.Range("J4").Value = "bobina"
$ws.Range("K4").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L4").Value = 8
$ws.Range("M4").Value = 70
$ws.Range("N4").Value = 4

$ws.Range("A5").Value = 250866
$ws.Range("B5").Value = "R3"
$ws.Range("C5").Value = 102
$ws.Range("D5").Value = 104.4081632653061
$ws.Range("E5").Value = "2025-04-10 07:18:00"
$ws.Range("F5").Value = "2025-04-10 09:00:00"
$ws.Range("G5").Value = "2025-04-10 09:00:00"
$ws.Range("H5").Value = "2025-04-10 10:44:24"
$ws.Range("I5").Value = 5116
$ws.Range("J5").Value = "bobina"
$ws.Range("K5").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L5").Value = 6
$ws.Range("M5").Value = 70
$ws.Range("N5").Value = 4

$ws.Range("A6").Value = 251550
$ws.Range("B6").Value = "R3"
$ws.Range("C6").Value = 50
$ws.Range("D6").Value = 727.5714285714286
$ws.Range("E6").Value = "2025-04-10 10:44:24"
$ws.Range("F6").Value = "2025-04-10 11:34:24"
$ws.Range("G6").Value = "2025-04-10 11:34:24"
$ws.Range("H6").Value = "2025-04-14 07:41:58"
$ws.Range("I6").Value = 35651
$ws.Range("J6").Value = "bobina"
$ws.Range("K6").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L6").Value = 3
$ws.Range("M6").Value = 70
$ws.Range("N6").Value = 4

$ws.Range("A7").Value = 251109
$ws.Range("B7").Value = "R6"
$ws.Range("C7").Value = 112
$ws.Range("D7").Value = 266.5915492957747
$ws.Range("E7").Value = "2025-04-10 13:25:00"
$ws.Range("F7").Value = "2025-04-11 07:17:00"
$ws.Range("G7").Value = "2025-04-11 07:17:00"
$ws.Range("H7").Value = "2025-04-11 11:43:35"
$ws.Range("I7").Value = 18928
$ws.Range("J7").Value = "bobina"
$ws.Range("K7").Value = "R6"
$ws.Range("L7").Value = 16
$ws.Range("M7").Value = 70
$ws.Range("N7").Value = 4

$ws.Range("A8").Value = 251088
$ws.Range("B8").Value = "R9"
$ws.Range("C8").Value = 35
$ws.Range("D8").Value = 89.6376811594203
$ws.Range("E8").Value = "2025-04-10 07:22:00"
$ws.Range("F8").Value = "2025-04-10 07:57:00"
$ws.Range("G8").Value = "2025-04-10 07:57:00"
$ws.Range("H8").Value = "2025-04-10 09:26:38"
$ws.Range("I8").Value = 6185
$ws.Range("J8").Value = "bobina"
$ws.Range("K8").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R10 ;R3 ;R6 ;R9"
$ws.Range("L8").Value = 3
$ws.Range("M8").Value = 70
$ws.Range("N8").Value = 4

$ws.Range("A9").Value = 235572
$ws.Range("B9").Value = "R9"
$ws.Range("C9").Value = 35
$ws.Range("D9").Value = 144.3188405797102
$ws.Range("E9").Value = "2025-04-10 09:26:38"
$ws.Range("F9").Value = "2025-04-10 10:01:38"
$ws.Range("G9").Value = "2025-04-10 10:01:38"
$ws.Range("H9").Value = "2025-04-10 12:25:57"
$ws.Range("I9").Value = 9958
$ws.Range("J9").Value = "bobina"
$ws.Range("K9").Value = "BIMEC 2 ;BIMEC 4 ;BIMEC 5 ;CASON ;R6 ;R9"
$ws.Range("L9").Value = 5
$ws.Range("M9").Value = 70
$ws.Range("N9").Value = 4
